# Logboek bijgewerkt met afgemaakte en niet afgemaakte onderdelen
#
# - Voeg een "Datum" kolom toe aan de tabel (Tabel2), met een datumnotatie.
# - Vul de datums in voor de taken die al een status (Ja) hebben.
# - Werk de status van twee taken bij van "Nee" naar "Ja".
# - Verplaats de actieve selectie naar J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Nieuwe kolom "Datum" toevoegen aan de tabel ---------------------
$lo = $ws.ListObjects.Item(1)
[void]$lo.ListColumns.Add()
$ws.Range("G3").Value = "Datum"

# --- 2. Statussen bijwerken van "Nee" naar "Ja" --------------------------
$ws.Range("F6").Value = "Ja"
$ws.Range("F7").Value = "Ja"

# --- 3. Datums invullen voor de afgeronde taken --------------------------
function Set-DatumCel($addr, [int]$serial) {
    $ws.Range($addr).Value2 = $serial
    $ws.Range($addr).NumberFormat = "d-mmm"
}

Set-DatumCel "G4"  43360   # 17 sep 2018
Set-DatumCel "G6"  43361   # 18 sep 2018
Set-DatumCel "G7"  43362   # 19 sep 2018
Set-DatumCel "G11" 43361   # 18 sep 2018
Set-DatumCel "G12" 43361   # 18 sep 2018

# --- 4. Selectie verplaatsen naar J11 ------------------------------------
[void]$ws.Range("J11").Select()
